$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "41.606.98"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.217.06"
$ws.Range("E3").Value = "  -1.80%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.46%  "

# Row 5 - BNB
Set-TextValue "D5" "230.31"
$ws.Range("E5").Value = "  -0.92%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -2.76%  "

# Row 7 - Solana
Set-TextValue "D7" "59.71"
$ws.Range("E7").Value = "  -7.21%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.14%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.402"
$ws.Range("E9").Value = "  -2.34%  "

# Row 10 - OKB
Set-TextValue "D10" "57.69"
$ws.Range("E10").Value = "  -2.93%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0893"
$ws.Range("E11").Value = "  -1.03%  "

# Row 12 - TRON
Set-TextValue "D12" "0.103"
$ws.Range("E12").Value = "  -1.32%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "2.551.68"
$ws.Range("E13").Value = "  -1.34%  "

# Row 14 - Chainlink
Set-TextValue "D14" "15.43"
$ws.Range("E14").Value = "  -5.93%  "

# Row 15 - Avalanche
Set-TextValue "D15" "22.40"
$ws.Range("E15").Value = "  -0.99%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.60"

# Row 17 - Polygon
Set-TextValue "D17" "0.797"
$ws.Range("E17").Value = "  -4.89%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.236.08"
$ws.Range("E18").Value = "  -1.12%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "41.596.60"
$ws.Range("E19").Value = "  +0.57%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -1.35%  "

# Row 21 - Litecoin
Set-TextValue "D21" "72.19"
$ws.Range("E21").Value = "  -2.37%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -2.22%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "244.89"
$ws.Range("E23").Value = "  -2.68%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.13%  "

# Row 25 / 26 swap: Toncoin <-> PancakeSwap
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D25" "2.36"
$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D26" "2.30"
$ws.Range("E26").Value = "  -5.07%  "

# Row 27 - Cosmos
Set-TextValue "D27" "9.72"
$ws.Range("E27").Value = "  -1.62%  "

# Row 28 - Monero
Set-TextValue "D28" "168.88"
$ws.Range("E28").Value = "  -2.76%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  -4.28%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "19.78"
$ws.Range("E30").Value = "  -3.62%  "

# Row 31 - ImmutableX
Set-TextValue "D31" "1.41"
$ws.Range("E31").Value = "  -4.51%  "

# Row 32 - WEMIXToken
Set-TextValue "D32" "2.52"
$ws.Range("E32").Value = "  -10.42%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  -2.86%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -1.84%  "

# Row 35 - Filecoin
Set-TextValue "D35" "4.66"
$ws.Range("E35").Value = "  -2.29%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  +1.98%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "2.37"
$ws.Range("E38").Value = "  -3.79%  "

# Row 39 - RenderToken
$ws.Range("E39").Value = "  -7.64%  "

# Row 40 - BinanceUSD
$ws.Range("E40").Value = "  +0.93%  "

# Row 41 - TerraClassic
Set-TextValue "D41" "0.000234"
$ws.Range("E41").Value = "  -12.55%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  -1.00%  "

# Row 43 - FraxShare
Set-TextValue "D43" "8.62"
$ws.Range("E43").Value = "  -3.05%  "

# Row 44 - Cronos
Set-TextValue "D44" "0.0963"
$ws.Range("E44").Value = "  +1.88%  "

# Row 45 / 46 swap: Aave <-> TrustWalletToken
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D45" "1.21"
$ws.Range("E45").Value = "  -2.27%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "97.63"
$ws.Range("E46").Value = "  -5.28%  "

# Row 47 - FTXToken
Set-TextValue "D47" "4.37"
$ws.Range("E47").Value = "  -10.51%  "

# Row 48 - Maker
Set-TextValue "D48" "1.469.30"
$ws.Range("E48").Value = "  -2.93%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "16.53"
$ws.Range("E49").Value = "  -8.04%  "

# Row 50 - HuobiToken
Set-TextValue "D50" "2.75"
$ws.Range("E50").Value = "  -1.74%  "

# Row 51 - ARBITRUM
$ws.Range("E51").Value = "  -4.80%  "
